$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two oldest years (2008年, 2009年) - this shifts every
# subsequent row up by two positions (old row4/2010年 becomes row2, etc.)
$ws.Rows("2:3").Delete()

# Append a new row for 2021年 at the bottom of the table (row 13)
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = 666
$ws.Range("D13").Value = 7290
$ws.Range("E13").Value = 1813
$ws.Range("F13").Value = 28.94
$ws.Range("G13").Value = 290.51
$ws.Range("H13").Value = 4.41
$ws.Range("I13").Value = 165.79
$ws.Range("J13").Value = 11.6182
$ws.Range("K13").Value = 2598
$ws.Range("L13").Value = 12.02984
$ws.Range("M13").Value = 31
$ws.Range("N13").Value = 284
$ws.Range("O13").Value = 1263
$ws.Range("P13").Value = ""
$ws.Range("Q13").Value = 99
$ws.Range("R13").Value = 659
$ws.Range("S13").Value = 10

# Match the formatting used by the rest of column A (bold, centered,
# thin box border) on the newly-added year label cell by copying the
# style already used by the cell directly above it.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
